# "Generate Report for Handback"
#
# The localization-status report is re-generated after handback:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     on the per-language sheets (zh-cn, de-de).
#   - Two new columns get populated for each file row: "Latest Target File"
#     (mirrors the source .md link) and "Latest Handback File" (mirrors the
#     handoff .xlf link) - both rendered as hyperlinks matching the existing
#     link styling (underline, same link color).
#   - "Latest Handback DateTime" moves from the zero-date sentinel to the
#     real handback timestamp (distinct per language sheet).

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2        # xlUnderlineStyleSingle
$hyperlinkColor = 15570276     # RGB(0x64,0x95,0xED) == style's FF6495ED font color

function Set-LinkCell {
    param($ws, [string]$cellRef, [string]$text, [string]$url)

    $ws.Range($cellRef).Value = $text
    # Hyperlinks.Add applies Excel's builtin theme-colored "Hyperlink" cell
    # style, clobbering any font formatting applied beforehand - so (re)apply
    # the underline + link color *after* adding the hyperlink to match the
    # workbook's existing custom hyperlink look (FF6495ED).
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $text)
    $ws.Range($cellRef).Font.Underline = $hyperlinkUnderline
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

Set-LinkCell $wsZh "F2" "56fba474-6806-4a0b-bcc0-d0124cb48364.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/18121388777d9e32d0db3a7b09c79a2f702284a3/e2e/56fba474-6806-4a0b-bcc0-d0124cb48364.md"
Set-LinkCell $wsZh "G2" "56fba474-6806-4a0b-bcc0-d0124cb48364.7f9261651327e5305f967cdabe6e4e8c6f094a33.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20c7c81b36cdd202bc9cb9ed5976b1c9c5388113/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/56fba474-6806-4a0b-bcc0-d0124cb48364.7f9261651327e5305f967cdabe6e4e8c6f094a33.zh-cn.xlf"

Set-LinkCell $wsZh "F3" "90e18c39-a61d-45bf-931b-8b10e7375373.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/18121388777d9e32d0db3a7b09c79a2f702284a3/e2e/90e18c39-a61d-45bf-931b-8b10e7375373.md"
Set-LinkCell $wsZh "G3" "90e18c39-a61d-45bf-931b-8b10e7375373.bb423061990ec528c344f1616d98afb155ff1842.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20c7c81b36cdd202bc9cb9ed5976b1c9c5388113/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/90e18c39-a61d-45bf-931b-8b10e7375373.bb423061990ec528c344f1616d98afb155ff1842.zh-cn.xlf"

$wsZh.Range("H2").Value = "2016-03-20 12:51:05"
$wsZh.Range("H3").Value = "2016-03-20 12:51:05"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

Set-LinkCell $wsDe "F2" "56fba474-6806-4a0b-bcc0-d0124cb48364.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/18121388777d9e32d0db3a7b09c79a2f702284a3/e2e/56fba474-6806-4a0b-bcc0-d0124cb48364.md"
Set-LinkCell $wsDe "G2" "56fba474-6806-4a0b-bcc0-d0124cb48364.7f9261651327e5305f967cdabe6e4e8c6f094a33.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7ed0318b44625e4affc336f637f4f2b32675edd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/56fba474-6806-4a0b-bcc0-d0124cb48364.7f9261651327e5305f967cdabe6e4e8c6f094a33.de-de.xlf"

Set-LinkCell $wsDe "F3" "90e18c39-a61d-45bf-931b-8b10e7375373.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/18121388777d9e32d0db3a7b09c79a2f702284a3/e2e/90e18c39-a61d-45bf-931b-8b10e7375373.md"
Set-LinkCell $wsDe "G3" "90e18c39-a61d-45bf-931b-8b10e7375373.bb423061990ec528c344f1616d98afb155ff1842.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7ed0318b44625e4affc336f637f4f2b32675edd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/90e18c39-a61d-45bf-931b-8b10e7375373.bb423061990ec528c344f1616d98afb155ff1842.de-de.xlf"

$wsDe.Range("H2").Value = "2016-03-20 12:51:10"
$wsDe.Range("H3").Value = "2016-03-20 12:51:10"
